$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cylinder-group summary table rows were reordered: the "6 cylinders /
#     carb 0" row moves up to row 7, the "4 cylinders" row moves down to
#     row 8, and the second "6 cylinders" row (row 9) now shows its own
#     "6" label instead of being merged with row 8. ------------------------

# Row 8's "E" cell was merged with row 9 (E8:E9) so that a single "6"
# spanned both rows; split them apart so each row can carry its own value.
$ws.Range("E8:E9").UnMerge()

# Give the (now unmerged) E8/E9 cells the same look as the other data
# cells in the "Cylinder" column (E7, E10) instead of the special
# vertical-centred merge style they used to have.
$ws.Range("E7").Copy()
$ws.Range("E8").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("E9").PasteSpecial(-4122)   # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Row 7 now holds the values that used to be in row 8
# (Cylinder=6, Engine/carb=0, N=2, Mean HP=110, SD HP=0, Mean Wt=2.7475, SD Wt=0.1803122292025695)
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 2
$ws.Range("H7").Value = 110
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 2.7475
$ws.Range("K7").Value = 0.1803122292025695

# Row 8 now holds the values that used to be in row 7
# (Cylinder=4, Engine/carb=1, N=3, Mean HP=83.33333333333333, SD HP=18.50225211517056,
#  Mean Wt=2.886666666666667, SD Wt=0.4911551010967242)
$ws.Range("E8").Value = 4
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3
$ws.Range("H8").Value = 83.33333333333333
$ws.Range("I8").Value = 18.50225211517056
$ws.Range("J8").Value = 2.886666666666667
$ws.Range("K8").Value = 0.4911551010967242

# Row 9 keeps its original F..K values; only its Cylinder label is now
# filled in explicitly (it used to be blank because it was merged into E8).
$ws.Range("E9").Value = 6
